# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# - Turn the data range into a real Excel Table (ListObject)
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = "_old" -> "_FV2210", L1:U1 = "_new" -> "_FV2304") ---
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2210"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2304"
}
# K1 ("diff") is unchanged.

# --- 2. Convert the used range into an Excel Table ---
$dataRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the top (header) row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
